$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 949.75
$ws.Range("J19").Value = 933
$ws.Range("L19").Value = 933
$ws.Range("N19").Value = -1283
$ws.Range("H64").Value = 4416.5
$ws.Range("I64").Value = 4416.5
$ws.Range("J64").Value = 0
$ws.Range("K64").Value = 4416.5
$ws.Range("L64").Value = 0
$ws.Range("M64").Value = -4168.5
$ws.Range("H67").Value = 4416.5
$ws.Range("I67").Value = 4416.5
$ws.Range("J67").Value = 0
$ws.Range("K67").Value = 4416.5
$ws.Range("L67").Value = 0
$ws.Range("M67").Value = -3558.5
$ws.Range("H86").Value = 3946.6667
$ws.Range("I86").Value = 3923.25
$ws.Range("J86").Value = 3958.375
$ws.Range("K86").Value = 3923.25
$ws.Range("L86").Value = 3958.375
$ws.Range("M86").Value = -2800.25
$ws.Range("N86").Value = -6204.375
$ws.Range("H89").Value = 3946.6667
$ws.Range("I89").Value = 3923.25
$ws.Range("J89").Value = 3958.375
$ws.Range("K89").Value = 19616.25
$ws.Range("L89").Value = 19791.875
$ws.Range("M89").Value = -14000.25
$ws.Range("N89").Value = -31023.875
$ws.Range("H106").Value = 2003.75
$ws.Range("I106").Value = 2003.75
$ws.Range("K106").Value = 2003.75
$ws.Range("M106").Value = -1372.75
$ws.Range("H107").Value = 1023.1
$ws.Range("I107").Value = 1066.5264
$ws.Range("K107").Value = 1066.5264
$ws.Range("M107").Value = 853.4736
$ws.Range("H113").Value = 4880.8
$ws.Range("I113").Value = 4452
$ws.Range("J113").Value = 5166.6665
$ws.Range("K113").Value = 4452
$ws.Range("L113").Value = 5166.6665
$ws.Range("M113").Value = -1198
$ws.Range("N113").Value = -11674.6665
$ws.Range("N64").ClearContents()
$ws.Range("N67").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 12108.216
$ws.Range("I32").Value = 11062.611
$ws.Range("K32").Value = 11062.611
$ws.Range("M32").Value = -10775.611
$ws.Range("H61").Value = 3489.25
$ws.Range("I61").Value = 3489.25
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 3489.25
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -3277.25
$ws.Range("H102").Value = 2363.5
$ws.Range("I102").Value = 2396.5454
$ws.Range("K102").Value = 2396.5454
$ws.Range("M102").Value = -774.5454
$ws.Range("H122").Value = 2850.6
$ws.Range("I122").Value = 2563.25
$ws.Range("J122").Value = 4000
$ws.Range("K122").Value = 7689.75
$ws.Range("L122").Value = 12000
$ws.Range("M122").Value = -5239.75
$ws.Range("N122").Value = -16900
$ws.Range("H136").Value = 3489.25
$ws.Range("I136").Value = 3489.25
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 10467.75
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -7917.75
$ws.Range("N61").ClearContents()
$ws.Range("N136").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 2049.375
$ws.Range("I107").Value = 1899.3334
$ws.Range("J107").Value = 2499.5
$ws.Range("K107").Value = 1899.3334
$ws.Range("L107").Value = 2499.5
$ws.Range("M107").Value = 20.66660000000002
$ws.Range("N107").Value = -6339.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 3313
$ws.Range("I16").Value = 2500
$ws.Range("K16").Value = 2500
$ws.Range("M16").Value = -2213
$ws.Range("H31").Value = 2169.1875
$ws.Range("I31").Value = 1739.7273
$ws.Range("J31").Value = 3114
$ws.Range("K31").Value = 1739.7273
$ws.Range("L31").Value = 3114
$ws.Range("M31").Value = -1444.7273
$ws.Range("N31").Value = -3704
$ws.Range("H34").Value = 2169.1875
$ws.Range("I34").Value = 1739.7273
$ws.Range("J34").Value = 3114
$ws.Range("K34").Value = 1739.7273
$ws.Range("L34").Value = 3114
$ws.Range("M34").Value = -1537.7273
$ws.Range("N34").Value = -3518
$ws.Range("H58").Value = 4567.4287
$ws.Range("I58").Value = 2693
$ws.Range("K58").Value = 2693
$ws.Range("M58").Value = -2490
$ws.Range("H113").Value = 3313
$ws.Range("I113").Value = 2500
$ws.Range("K113").Value = 2500
$ws.Range("M113").Value = -330
$ws.Range("H136").Value = 4567.4287
$ws.Range("I136").Value = 2693
$ws.Range("K136").Value = 8079
$ws.Range("M136").Value = -5529
$ws.Range("H141").Value = 227993.42
$ws.Range("J141").Value = 227993.42
$ws.Range("L141").Value = 227993.42
$ws.Range("N141").Value = -238353.42

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H37").Value = 149849.5
$ws.Range("J37").Value = 149849.5
$ws.Range("L37").Value = 449548.5
$ws.Range("N37").Value = -449772.5
$ws.Range("H136").Value = 2367.3333
$ws.Range("I136").Value = 2005
$ws.Range("J136").Value = 2423.077
$ws.Range("K136").Value = 6015
$ws.Range("L136").Value = 7269.231000000001
$ws.Range("M136").Value = -915
$ws.Range("N136").Value = -17469.231
$ws.Range("H138").Value = 3150
$ws.Range("I138").Value = 1050
$ws.Range("K138").Value = 3150
$ws.Range("M138").Value = 1990
$ws.Range("H139").Value = 1380.2
$ws.Range("I139").Value = 1380.2
$ws.Range("K139").Value = 4140.6
$ws.Range("M139").Value = 999.3999999999996
$ws.Range("H140").Value = 3074.875
$ws.Range("I140").Value = 3074.875
$ws.Range("K140").Value = 9224.625
$ws.Range("M140").Value = -4044.625

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 4144.6
$ws.Range("I132").Value = 3690
$ws.Range("J132").Value = 4599.2
$ws.Range("K132").Value = 11070
$ws.Range("L132").Value = 13797.6
$ws.Range("M132").Value = -8540
$ws.Range("N132").Value = -18857.6

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 5300
$ws.Range("I40").Value = 4950
$ws.Range("K40").Value = 4950
$ws.Range("M40").Value = -4814
$ws.Range("H68").Value = 2722.65
$ws.Range("I68").Value = 2750.158
$ws.Range("K68").Value = 2750.158
$ws.Range("M68").Value = -2001.158
$ws.Range("H71").Value = 2722.65
$ws.Range("I71").Value = 2750.158
$ws.Range("K71").Value = 13750.79
$ws.Range("M71").Value = -10006.79
$ws.Range("H82").Value = 3036
$ws.Range("J82").Value = 3850
$ws.Range("L82").Value = 3850
$ws.Range("N82").Value = -4572
$ws.Range("H85").Value = 3036
$ws.Range("J85").Value = 3850
$ws.Range("L85").Value = 3850
$ws.Range("N85").Value = -6346
$ws.Range("H122").Value = 3810.5
$ws.Range("I122").Value = 3772.6
$ws.Range("K122").Value = 11317.8
$ws.Range("M122").Value = -8867.799999999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 2624.0908
$ws.Range("I81").Value = 1996.1111
$ws.Range("J81").Value = 5450
$ws.Range("K81").Value = 3992.2222
$ws.Range("L81").Value = 10900
$ws.Range("M81").Value = -2931.2222
$ws.Range("N81").Value = -13022
$ws.Range("H84").Value = 2624.0908
$ws.Range("I84").Value = 1996.1111
$ws.Range("J84").Value = 5450
$ws.Range("K84").Value = 19961.111
$ws.Range("L84").Value = 54500
$ws.Range("M84").Value = -14657.111
$ws.Range("N84").Value = -65108
$ws.Range("H122").Value = 1214.7778
$ws.Range("J122").Value = 1337.5
$ws.Range("L122").Value = 4012.5
$ws.Range("N122").Value = -8912.5
$ws.Range("H126").Value = 1700
$ws.Range("J126").Value = 1700
$ws.Range("L126").Value = 5100
$ws.Range("N126").Value = -10040
$ws.Range("H132").Value = 7199.3335
$ws.Range("I132").Value = 2733
$ws.Range("K132").Value = 8199
$ws.Range("M132").Value = -5669
